# Append 9 new data rows (rows 11-19) to the Andre Russell stats sheet,
# extending the table from A1:K10 to A1:K19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0
$batsman = "Andre Russell" + $nbsp

# Each entry: venue, date, result, ownTeam, oppTeam, batsman, runs, balls, 4s, 6s, sr
$rows = @(
  @(" Abu Dhabi",    " September 23 2020", "Mumbai won by 49 runs",                              "Kolkata Knight Riders", "Mumbai Indians",               $batsman, "11", "11", "2", "0", "100.00"),
  @(" Abu Dhabi",    " October 18 2020",   "Match tied (KKR won the one-over eliminator)",        "Kolkata Knight Riders", "Sunrisers Hyderabad",          $batsman, "9",  "11", "1", "0", "81.81"),
  @(" Dubai (DSC)",  " September 30 2020", "KKR won by 37 runs",                                  "Kolkata Knight Riders", "Rajasthan Royals",             $batsman, "24", "14", "0", "3", "171.42"),
  @(" Abu Dhabi",    " October 16 2020",   "Mumbai won by 8 wickets (with 19 balls remaining)",   "Kolkata Knight Riders", "Mumbai Indians",               $batsman, "12", "9",  "1", "1", "133.33"),
  @(" Dubai (DSC)",  " November 01 2020",  "KKR won by 60 runs",                                  "Kolkata Knight Riders", "Rajasthan Royals",             $batsman, "25", "11", "1", "3", "227.27"),
  @(" Abu Dhabi",    " October 10 2020",   "KKR won by 2 runs",                                   "Kolkata Knight Riders", "Kings XI Punjab",              $batsman, "5",  "3",  "1", "0", "166.66"),
  @(" Sharjah",      " October 03 2020",   "Capitals won by 18 runs",                             "Kolkata Knight Riders", "Delhi Capitals",               $batsman, "13", "8",  "1", "1", "162.50"),
  @(" Sharjah",      " October 12 2020",   "RCB won by 82 runs",                                  "Kolkata Knight Riders", "Royal Challengers Bangalore",  $batsman, "16", "10", "2", "1", "160.00"),
  @(" Abu Dhabi",    " October 07 2020",   "KKR won by 10 runs",                                  "Kolkata Knight Riders", "Chennai Super Kings",          $batsman, "2",  "4",  "0", "0", "50.00")
)

$startRow = 11
$endRow = $startRow + $rows.Length - 1

# Force the newly written numeric-looking values (columns G:K) to be stored
# as text, matching the existing rows (which are all text cells).
$ws.Range("A" + $startRow + ":K" + $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    $ws.Range("A" + $r).Value = $vals[0]
    $ws.Range("B" + $r).Value = $vals[1]
    $ws.Range("C" + $r).Value = $vals[2]
    $ws.Range("D" + $r).Value = $vals[3]
    $ws.Range("E" + $r).Value = $vals[4]
    $ws.Range("F" + $r).Value = $vals[5]
    $ws.Range("G" + $r).Value = $vals[6]
    $ws.Range("H" + $r).Value = $vals[7]
    $ws.Range("I" + $r).Value = $vals[8]
    $ws.Range("J" + $r).Value = $vals[9]
    $ws.Range("K" + $r).Value = $vals[10]
}
